# Update the GWL10 / OEKS15 lookup table: each model's re-calculated
# "Period" (column B) and "Gap (Years)" (column C) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ B = "2045-2064"; C = 44 },
    @{ B = "2034-2053"; C = 33 },
    @{ B = "2034-2053"; C = 33 },
    @{ B = "2036-2055"; C = 35 },
    @{ B = "2021-2040"; C = 20 },
    @{ B = "2027-2046"; C = 26 },
    @{ B = "n/a";       C = "n/a" },
    @{ B = "n/a";       C = "n/a" },
    @{ B = "2029-2048"; C = 28 },
    @{ B = "n/a";       C = "n/a" },
    @{ B = "2041-2060"; C = 40 },
    @{ B = "2020-2039"; C = 19 },
    @{ B = "2020-2039"; C = 19 },
    @{ B = "2026-2045"; C = 25 },
    @{ B = "2030-2049"; C = 29 },
    @{ B = "2016-2035"; C = 15 },
    @{ B = "2010-2029"; C = 9 },
    @{ B = "2017-2036"; C = 16 },
    @{ B = "2021-2040"; C = 20 },
    @{ B = "2020-2039"; C = 19 },
    @{ B = "2001-2020"; C = 0 },
    @{ B = "2006-2025"; C = 5 },
    @{ B = "2005-2024"; C = 4 },
    @{ B = "2008-2027"; C = 7 },
    @{ B = "2006-2025"; C = 5 },
    @{ B = "2005-2024"; C = 4 },
    @{ B = "2005-2024"; C = 4 },
    @{ B = "2029-2048"; C = 28 },
    @{ B = "2019-2038"; C = 18 },
    @{ B = "n/a";       C = "n/a" },
    @{ B = "2029-2048"; C = 28 },
    @{ B = "2022-2041"; C = 21 },
    @{ B = "2025-2044"; C = 24 },
    @{ B = "2025-2044"; C = 24 }
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i].B
    $ws.Cells.Item($row, 3).Value = $data[$i].C
}

$ws.Range("F23").Select()
